# "Moving the data to the right file"
# A new qPCR replicate column (D1) is inserted after the existing D2 column (K).
# The previous K-column values shift right into the new L column, and K is
# populated with a new set of replicate values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the inserted column
$ws.Range("L1").Value = "D1"

# Copy the number formatting/style from column K (data rows) onto column L
# before writing values, so the new cells match the existing "Normal 2" style.
$ws.Range("K2:K4").Copy() | Out-Null
$ws.Range("L2:L4").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Shift the old K column values into the new L column
$ws.Range("L2").Value = 32.272718109441101
$ws.Range("L3").Value = 32.354724267277099
$ws.Range("L4").Value = 32.752593145860203

# Write the new values into K
$ws.Range("K2").Value = 34.426564717522503
$ws.Range("K3").Value = 35.173269496183302
$ws.Range("K4").Value = 35.4959667909053

# Update the active selection to match the saved workbook state
$ws.Range("B6").Select() | Out-Null
